# "Generate Report for Handoff"
# Appends a new handoff record (the file
# 2e74b109-60c5-4598-9f3a-b9c4069f4780oo....md) as row 3 on every sheet:
#   - Overview : File Name / Path And Name / Extension / Publish URL / zh-cn / de-de / Latest HO Xliff Generate Date
#   - zh-cn    : Source File Name ... Error Detail
#   - de-de    : Source File Name ... Error Detail

$wb = $excel.ActiveWorkbook

$commit = "733efbf282cf8b2edff632922488872460dc0636"
$fileName = "2e74b109-60c5-4598-9f3a-b9c4069f4780ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathAndName = "e2e\" + $fileName
$githubUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commit + "/e2e/" + $fileName

$zhXlf = "2e74b109-60c5-4598-9f3a-b9c4069f4780oooooooooooooooooooooooooooooooooooooooo.6c29a9e4676665808b29457443c55e9059d852b3.zh-cn.xlf"
$deXlf = "2e74b109-60c5-4598-9f3a-b9c4069f4780oooooooooooooooooooooooooooooooooooooooo.6c29a9e4676665808b29457443c55e9059d852b3.de-de.xlf"

$status = "Ready for handoff"
$zhHandoffDate = "2016-08-19 06:27:48"
$deHandoffDate = "2016-08-19 06:27:54"
$overviewDate = "2016-08-19 06:27:54"

# ---------------------------------------------------------------------------
# Overview sheet (7-column table)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathAndName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $overviewDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $githubUrl, "", "", $pathAndName) | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# zh-cn sheet (16-column table)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $fileName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $githubUrl, "", "", $fileName) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------------
# de-de sheet (16-column table)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $fileName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $githubUrl, "", "", $fileName) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 16.3
